# Insert a new weekly record at row 55 (Hortaliza, Terminal Hortofrutícola Agro
# Chillán - Pepino ensalada). This pushes the existing rows 55..122 down to
# 56..123 and adds a brand-new data row in the now-empty row 55.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(55).Insert()

$ws.Range("A55").Value = 7
$ws.Range("B55").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C55").Value = "Ñuble"
$ws.Range("D55").Value = 44413
$ws.Range("E55").Value = 16
$ws.Range("F55").Value = 100112043
$ws.Range("G55").Value = "Pepino ensalada"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 120
$ws.Range("K55").Value = 17000
$ws.Range("L55").Value = 18000
$ws.Range("M55").Value = 17500
$ws.Range("N55").Value = "$/caja 60 unidades"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 292
$ws.Range("Q55").Value = 60
$ws.Range("R55").Value = "Hortaliza"
